# Update "想去人数" (interest count) values in column F across sheets,
# reflecting the regenerated data output.

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 509
$ws1.Range("F4").Value = 509
$ws1.Range("F5").Value = 789
$ws1.Range("F6").Value = 187
$ws1.Range("F8").Value = 808
$ws1.Range("F10").Value = 614
$ws1.Range("F11").Value = 168
$ws1.Range("F15").Value = 105
$ws1.Range("F16").Value = 1535
$ws1.Range("F17").Value = 187
$ws1.Range("F20").Value = 56
$ws1.Range("F24").Value = 19
$ws1.Range("F25").Value = 207
$ws1.Range("F27").Value = 48
$ws1.Range("F28").Value = 1377
$ws1.Range("F29").Value = 129

# Sheet: 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 30
$ws2.Range("F6").Value = 7

# Sheet: 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 509
$ws4.Range("F5").Value = 509
$ws4.Range("F6").Value = 789
$ws4.Range("F7").Value = 187
$ws4.Range("F9").Value = 808
$ws4.Range("F11").Value = 30
$ws4.Range("F13").Value = 614
$ws4.Range("F15").Value = 168
$ws4.Range("F19").Value = 105
$ws4.Range("F20").Value = 1535
$ws4.Range("F22").Value = 187
$ws4.Range("F25").Value = 56
$ws4.Range("F27").Value = 7
$ws4.Range("F36").Value = 19
$ws4.Range("F37").Value = 207
$ws4.Range("F39").Value = 48
$ws4.Range("F40").Value = 1377
$ws4.Range("F41").Value = 129
